# Updates odds values in Sheet1 for rows 2 and 7, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62

# Row 7 updates
$ws.Range("G7").Value = 5.4
$ws.Range("H7").Value = 3.65
$ws.Range("I7").Value = 1.55
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 2.12
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.8
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.67
$ws.Range("U7").Value = 1.98
$ws.Range("V7").Value = 1.75
$ws.Range("W7").Value = 13
$ws.Range("X7").Value = 32
$ws.Range("Y7").Value = 18
$ws.Range("Z7").Value = 110
$ws.Range("AB7").Value = 70
$ws.Range("AD7").Value = 7.3
$ws.Range("AE7").Value = 18.5
$ws.Range("AG7").Value = 6.1
$ws.Range("AH7").Value = 6.9
$ws.Range("AJ7").Value = 11
$ws.Range("AK7").Value = 13
$ws.Range("AN7").Value = 7
$ws.Range("AP7").Value = 40
$ws.Range("AQ7").Value = 250
$ws.Range("AR7").Value = 300
$ws.Range("AT7").Value = 2.67
$ws.Range("AU7").Value = 8.25
$ws.Range("AV7").Value = 90
$ws.Range("AX7").Value = 7.6
$ws.Range("AY7").Value = 18.5
$ws.Range("AZ7").Value = 25
$ws.Range("BA7").Value = 60
$ws.Range("BB7").Value = 300
